$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row (351) is copied to each new row first so the new cells
# inherit identical (default) formatting without minting new style records,
# then the per-row values are written over the top.
$template = $ws.Range("A351:M351")

$template.Copy($ws.Range("A352:M352"))
$ws.Cells.Item(352, 1).Value = 351
$ws.Cells.Item(352, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(352, 3).Value = "3:05 PM"
$ws.Cells.Item(352, 4).Value = "SK759"
$ws.Cells.Item(352, 5).Value = "Copenhagen"
$ws.Cells.Item(352, 6).Value = "(CPH)"
$ws.Cells.Item(352, 7).Value = "SAS "
$ws.Cells.Item(352, 8).Value = "A20N"
$ws.Cells.Item(352, 9).Value = "(EI-SIE)"
$ws.Cells.Item(352, 10).Value = "3:00 PM"
$ws.Cells.Item(352, 12).Value = "0 hours, -5 minutes"

$template.Copy($ws.Range("A353:M353"))
$ws.Cells.Item(353, 1).Value = 352
$ws.Cells.Item(353, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(353, 3).Value = "3:50 PM"
$ws.Cells.Item(353, 4).Value = "W61732"
$ws.Cells.Item(353, 5).Value = "Stockholm"
$ws.Cells.Item(353, 6).Value = "(NYO)"
$ws.Cells.Item(353, 7).Value = "Wizz Air "
$ws.Cells.Item(353, 8).Value = "A321"
$ws.Cells.Item(353, 9).Value = "(HA-LTB)"
$ws.Cells.Item(353, 10).Value = "3:22 PM"
$ws.Cells.Item(353, 12).Value = "0 hours, -28 minutes"

$template.Copy($ws.Range("A354:M354"))
$ws.Cells.Item(354, 1).Value = 353
$ws.Cells.Item(354, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(354, 3).Value = "4:20 PM"
$ws.Cells.Item(354, 4).Value = "FR7272"
$ws.Cells.Item(354, 5).Value = "Vaxjo"
$ws.Cells.Item(354, 6).Value = "(VXO)"
$ws.Cells.Item(354, 7).Value = "Ryanair "
$ws.Cells.Item(354, 8).Value = "B738"
$ws.Cells.Item(354, 9).Value = "(SP-RKQ)"
$ws.Cells.Item(354, 10).Value = "4:16 PM"
$ws.Cells.Item(354, 12).Value = "0 hours, -4 minutes"

$template.Copy($ws.Range("A355:M355"))
$ws.Cells.Item(355, 1).Value = 354
$ws.Cells.Item(355, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(355, 3).Value = "4:25 PM"
$ws.Cells.Item(355, 4).Value = "LO3815"
$ws.Cells.Item(355, 5).Value = "Warsaw"
$ws.Cells.Item(355, 6).Value = "(WAW)"
$ws.Cells.Item(355, 7).Value = "LOT "
$ws.Cells.Item(355, 8).Value = "E195"
$ws.Cells.Item(355, 9).Value = "(SP-LNI)"
$ws.Cells.Item(355, 10).Value = "4:07 PM"
$ws.Cells.Item(355, 12).Value = "0 hours, -18 minutes"

$template.Copy($ws.Range("A356:M356"))
$ws.Cells.Item(356, 1).Value = 355
$ws.Cells.Item(356, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(356, 3).Value = "5:30 PM"
$ws.Cells.Item(356, 4).Value = "KL1921"
$ws.Cells.Item(356, 5).Value = "Amsterdam"
$ws.Cells.Item(356, 6).Value = "(AMS)"
$ws.Cells.Item(356, 7).Value = "KLM "
$ws.Cells.Item(356, 8).Value = "E190"
$ws.Cells.Item(356, 9).Value = "(PH-EZU)"
$ws.Cells.Item(356, 10).Value = "5:22 PM"
$ws.Cells.Item(356, 12).Value = "0 hours, -8 minutes"

$template.Copy($ws.Range("A357:M357"))
$ws.Cells.Item(357, 1).Value = 356
$ws.Cells.Item(357, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(357, 3).Value = "5:45 PM"
$ws.Cells.Item(357, 4).Value = "W61784"
$ws.Cells.Item(357, 5).Value = "Oslo"
$ws.Cells.Item(357, 6).Value = "(OSL)"
$ws.Cells.Item(357, 7).Value = "Wizz Air "
$ws.Cells.Item(357, 8).Value = "A320"
$ws.Cells.Item(357, 9).Value = "(HA-LYS)"
$ws.Cells.Item(357, 10).Value = "5:32 PM"
$ws.Cells.Item(357, 12).Value = "0 hours, -13 minutes"

$template.Copy($ws.Range("A358:M358"))
$ws.Cells.Item(358, 1).Value = 357
$ws.Cells.Item(358, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(358, 3).Value = "5:50 PM"
$ws.Cells.Item(358, 4).Value = "W61632"
$ws.Cells.Item(358, 5).Value = "Paris"
$ws.Cells.Item(358, 6).Value = "(BVA)"
$ws.Cells.Item(358, 7).Value = "Wizz Air "
$ws.Cells.Item(358, 8).Value = "A320"
$ws.Cells.Item(358, 9).Value = "(HA-LYO)"
$ws.Cells.Item(358, 10).Value = "5:18 PM"
$ws.Cells.Item(358, 12).Value = "0 hours, -32 minutes"

$template.Copy($ws.Range("A359:M359"))
$ws.Cells.Item(359, 1).Value = 358
$ws.Cells.Item(359, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(359, 3).Value = "5:55 PM"
$ws.Cells.Item(359, 4).Value = "W61746"
$ws.Cells.Item(359, 5).Value = "Bergen"
$ws.Cells.Item(359, 6).Value = "(BGO)"
$ws.Cells.Item(359, 7).Value = "Wizz Air "
$ws.Cells.Item(359, 8).Value = "A321"
$ws.Cells.Item(359, 9).Value = "(HA-LXL)"
$ws.Cells.Item(359, 10).Value = "5:35 PM"
$ws.Cells.Item(359, 12).Value = "0 hours, -20 minutes"

$template.Copy($ws.Range("A360:M360"))
$ws.Cells.Item(360, 1).Value = 359
$ws.Cells.Item(360, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(360, 3).Value = "6:10 PM"
$ws.Cells.Item(360, 4).Value = "FR239"
$ws.Cells.Item(360, 5).Value = "Krakow"
$ws.Cells.Item(360, 6).Value = "(KRK)"
$ws.Cells.Item(360, 7).Value = "Ryanair "
$ws.Cells.Item(360, 8).Value = "B738"
$ws.Cells.Item(360, 9).Value = "(SP-RSO)"
$ws.Cells.Item(360, 10).Value = "6:06 PM"
$ws.Cells.Item(360, 12).Value = "0 hours, -4 minutes"

$template.Copy($ws.Range("A361:M361"))
$ws.Cells.Item(361, 1).Value = 360
$ws.Cells.Item(361, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(361, 3).Value = "6:25 PM"
$ws.Cells.Item(361, 4).Value = "W61734"
$ws.Cells.Item(361, 5).Value = "Stockholm"
$ws.Cells.Item(361, 6).Value = "(NYO)"
$ws.Cells.Item(361, 7).Value = "Wizz Air "
$ws.Cells.Item(361, 8).Value = "A320"
$ws.Cells.Item(361, 9).Value = "(HA-LYM)"
$ws.Cells.Item(361, 10).Value = "6:04 PM"
$ws.Cells.Item(361, 12).Value = "0 hours, -21 minutes"

$template.Copy($ws.Range("A362:M362"))
$ws.Cells.Item(362, 1).Value = 361
$ws.Cells.Item(362, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(362, 3).Value = "6:30 PM"
$ws.Cells.Item(362, 4).Value = "FR3302"
$ws.Cells.Item(362, 5).Value = "Birmingham"
$ws.Cells.Item(362, 6).Value = "(BHX)"
$ws.Cells.Item(362, 7).Value = "Ryanair "
$ws.Cells.Item(362, 8).Value = "B738"
$ws.Cells.Item(362, 9).Value = "(SP-RSL)"
$ws.Cells.Item(362, 10).Value = "6:34 PM"
$ws.Cells.Item(362, 12).Value = "0 hours, 4 minutes"

$template.Copy($ws.Range("A363:M363"))
$ws.Cells.Item(363, 1).Value = 362
$ws.Cells.Item(363, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(363, 3).Value = "7:25 PM"
$ws.Cells.Item(363, 4).Value = "FR3272"
$ws.Cells.Item(363, 5).Value = "Dublin"
$ws.Cells.Item(363, 6).Value = "(DUB)"
$ws.Cells.Item(363, 7).Value = "Ryanair "
$ws.Cells.Item(363, 8).Value = "B738"
$ws.Cells.Item(363, 9).Value = "(SP-RSW)"
$ws.Cells.Item(363, 10).Value = "7:13 PM"
$ws.Cells.Item(363, 12).Value = "0 hours, -12 minutes"

$template.Copy($ws.Range("A364:M364"))
$ws.Cells.Item(364, 1).Value = 363
$ws.Cells.Item(364, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(364, 3).Value = "8:15 PM"
$ws.Cells.Item(364, 4).Value = "FR3286"
$ws.Cells.Item(364, 5).Value = "Leeds"
$ws.Cells.Item(364, 6).Value = "(LBA)"
$ws.Cells.Item(364, 7).Value = "Ryanair "
$ws.Cells.Item(364, 8).Value = "B738"
$ws.Cells.Item(364, 9).Value = "(SP-RKM)"
$ws.Cells.Item(364, 10).Value = "8:08 PM"
$ws.Cells.Item(364, 12).Value = "0 hours, -7 minutes"

$template.Copy($ws.Range("A365:M365"))
$ws.Cells.Item(365, 1).Value = 364
$ws.Cells.Item(365, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(365, 3).Value = "8:50 PM"
$ws.Cells.Item(365, 4).Value = "LO3825"
$ws.Cells.Item(365, 5).Value = "Warsaw"
$ws.Cells.Item(365, 6).Value = "(WAW)"
$ws.Cells.Item(365, 7).Value = "LOT (Retro Livery) "
$ws.Cells.Item(365, 8).Value = "E75S"
$ws.Cells.Item(365, 9).Value = "(SP-LIM)"
$ws.Cells.Item(365, 10).Value = "9:14 PM"
$ws.Cells.Item(365, 12).Value = "0 hours, 24 minutes"

$template.Copy($ws.Range("A366:M366"))
$ws.Cells.Item(366, 1).Value = 365
$ws.Cells.Item(366, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(366, 3).Value = "8:55 PM"
$ws.Cells.Item(366, 4).Value = "LH1646"
$ws.Cells.Item(366, 5).Value = "Munich"
$ws.Cells.Item(366, 6).Value = "(MUC)"
$ws.Cells.Item(366, 7).Value = "Lufthansa "
$ws.Cells.Item(366, 8).Value = "CRJ9"
$ws.Cells.Item(366, 9).Value = "(D-ACNM)"
$ws.Cells.Item(366, 10).Value = "8:48 PM"
$ws.Cells.Item(366, 12).Value = "0 hours, -7 minutes"

$template.Copy($ws.Range("A367:M367"))
$ws.Cells.Item(367, 1).Value = 366
$ws.Cells.Item(367, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(367, 3).Value = "8:58 PM"
$ws.Cells.Item(367, 4).Value = "UNKNOWN"
$ws.Cells.Item(367, 5).Value = "Tuzla"
$ws.Cells.Item(367, 6).Value = "(TZL)"
$ws.Cells.Item(367, 7).Value = "Wizz Air "
$ws.Cells.Item(367, 8).Value = "A320"
$ws.Cells.Item(367, 9).Value = "(HA-LYH)"
$ws.Cells.Item(367, 10).Value = "8:58 PM"
$ws.Cells.Item(367, 12).Value = "0 hours, 0 minutes"

$template.Copy($ws.Range("A368:M368"))
$ws.Cells.Item(368, 1).Value = 367
$ws.Cells.Item(368, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(368, 3).Value = "9:35 PM"
$ws.Cells.Item(368, 4).Value = "FR545"
$ws.Cells.Item(368, 5).Value = "London"
$ws.Cells.Item(368, 6).Value = "(STN)"
$ws.Cells.Item(368, 7).Value = "Ryanair "
$ws.Cells.Item(368, 8).Value = "B738"
$ws.Cells.Item(368, 9).Value = "(SP-RKQ)"
$ws.Cells.Item(368, 10).Value = "10:04 PM"
$ws.Cells.Item(368, 12).Value = "0 hours, 29 minutes"

$template.Copy($ws.Range("A369:M369"))
$ws.Cells.Item(369, 1).Value = 368
$ws.Cells.Item(369, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(369, 3).Value = "10:20 PM"
$ws.Cells.Item(369, 4).Value = "FR8509"
$ws.Cells.Item(369, 5).Value = "Oslo"
$ws.Cells.Item(369, 6).Value = "(TRF)"
$ws.Cells.Item(369, 7).Value = "Ryanair "
$ws.Cells.Item(369, 8).Value = "B738"
$ws.Cells.Item(369, 9).Value = "(SP-RSL)"
$ws.Cells.Item(369, 10).Value = "10:43 PM"
$ws.Cells.Item(369, 12).Value = "0 hours, 23 minutes"

$template.Copy($ws.Range("A370:M370"))
$ws.Cells.Item(370, 1).Value = 369
$ws.Cells.Item(370, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(370, 3).Value = "10:25 PM"
$ws.Cells.Item(370, 4).Value = "FR6112"
$ws.Cells.Item(370, 5).Value = "Lublin"
$ws.Cells.Item(370, 6).Value = "(LUZ)"
$ws.Cells.Item(370, 7).Value = "Ryanair "
$ws.Cells.Item(370, 8).Value = "B738"
$ws.Cells.Item(370, 9).Value = "(SP-RSW)"
$ws.Cells.Item(370, 10).Value = "10:33 PM"
$ws.Cells.Item(370, 12).Value = "0 hours, 8 minutes"

$template.Copy($ws.Range("A371:M371"))
$ws.Cells.Item(371, 1).Value = 370
$ws.Cells.Item(371, 2).Value = "Friday, Jan 13"
$ws.Cells.Item(371, 3).Value = "11:00 PM"
$ws.Cells.Item(371, 4).Value = "W61762"
$ws.Cells.Item(371, 5).Value = "Trondheim"
$ws.Cells.Item(371, 6).Value = "(TRD)"
$ws.Cells.Item(371, 7).Value = "Wizz Air "
$ws.Cells.Item(371, 8).Value = "A320"
$ws.Cells.Item(371, 9).Value = "(HA-LYO)"
$ws.Cells.Item(371, 10).Value = "10:37 PM"
$ws.Cells.Item(371, 12).Value = "0 hours, -23 minutes"
